$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orakeltjeneste")

# Oppdatert kursinfo for V26: fristen for innleveringssettet er flyttet
# fra 25.04 til 17.04.
$ws.Range("C2").Value = "Innleveringssett har frist 17.04"

# Reflect the author's new active selection after the edit.
$ws.Range("C3").Select()
